$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 229, shifting existing rows 229:239 down to 230:240.
$ws.Rows("229:229").Insert()

# Populate the newly inserted row 229 with the new weekly record.
$ws.Range("A229").Value = 9
$ws.Range("B229").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C229").Value = "Metropolitana"
$ws.Range("D229").Value = 44610
$ws.Range("E229").Value = 13
$ws.Range("F229").Value = 100112043
$ws.Range("G229").Value = "Pepino ensalada"
$ws.Range("H229").Value = "Sin especificar"
$ws.Range("I229").Value = "Primera"
$ws.Range("J229").Value = 106
$ws.Range("K229").Value = 10000
$ws.Range("L229").Value = 11000
$ws.Range("M229").Value = 10500
$ws.Range("N229").Value = "$/caja 70 unidades"
$ws.Range("O229").Value = "Región del Maule"
$ws.Range("P229").Value = 150
$ws.Range("Q229").Value = 70
$ws.Range("R229").Value = "Hortaliza"
